$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row copied into I3:L3
$ws.Range("I3").Value = "1 mg/ml"
$ws.Range("J3").Value = "0.1 mg/ml"
$ws.Range("K3").Value = "0.01 mg/ml"
$ws.Range("L3").Value = "Control"

# Averages in row 4
$ws.Range("I4").Formula = "=AVERAGE(A2:A22)"
$ws.Range("J4").Formula = "=AVERAGE(B2:B22)"
$ws.Range("K4").Formula = "=AVERAGE(C2:C22)"
$ws.Range("L4").Formula = "=AVERAGE(D2:D22)"

# Std devs in row 5
$ws.Range("I5").Formula = "=STDEV.S(A2:A22)"
$ws.Range("J5").Formula = "=STDEV.S(B2:B22)"
$ws.Range("K5").Formula = "=STDEV.S(C2:C22)"
$ws.Range("L5").Formula = "=STDEV.S(D2:D22)"

# Header row copied into I8:L8
$ws.Range("I8").Value = "1 mg/ml"
$ws.Range("J8").Value = "0.1 mg/ml"
$ws.Range("K8").Value = "0.01 mg/ml"
$ws.Range("L8").Value = "Control"

# Normalized row 9 (percent style)
$ws.Range("I9").Formula = '=I4/$L$4'
$ws.Range("J9").Formula = '=J4/$L$4'
$ws.Range("K9").Formula = '=K4/$L$4'
$ws.Range("L9").Formula = '=L4/$L$4'

# Normalized row 10 (percent style)
$ws.Range("I10").Formula = '=I5/$L$4'
$ws.Range("J10").Formula = '=J5/$L$4'
$ws.Range("K10").Formula = '=K5/$L$4'
$ws.Range("L10").Formula = '=L5/$L$4'

$ws.Range("I9:L10").Style = "Percent"

$ws.Range("M14").Select()
